# Fixed issue with setup of FCF in imports
# Adds two new columns (S: CF1, T: CF2) to the investor KYC import sheet,
# with row values 1/A, 2/B, 3/C, 4/D for the four data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1) - new columns S1 and T1
$ws.Range("S1").Value = "CF1"
$ws.Range("T1").Value = "CF2"

# Data rows 2-5: S column gets sequential numbers, T column gets letters
$ws.Cells.Item(2, 19).Value = 1
$ws.Cells.Item(2, 20).Value = "A"

$ws.Cells.Item(3, 19).Value = 2
$ws.Cells.Item(3, 20).Value = "B"

$ws.Cells.Item(4, 19).Value = 3
$ws.Cells.Item(4, 20).Value = "C"

$ws.Cells.Item(5, 19).Value = 4
$ws.Cells.Item(5, 20).Value = "D"

# Move the active selection/viewport to reflect the extended data range,
# as seen in the saved workbook (scrolled so column D is leftmost,
# active cell T6).
$ws.Range("T6").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
